# Update crypto price/volume data (and row reshuffle for some coins)
# as scraped on Wed Jun 21 23:30:01 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "250.19",
    # "30.105.35") are not auto-coerced into Excel numbers/dates,
    # matching the inlineStr text cells in the source file.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '30.105.35'
Set-TextCell 'E2' '  +6.69%  '
# Row 3
Set-TextCell 'D3' '1.894.51'
Set-TextCell 'E3' '  +5.98%  '
# Row 4
Set-TextCell 'E4' '  -0.11%  '
# Row 5
Set-TextCell 'D5' '250.19'
Set-TextCell 'E5' '  +1.71%  '
# Row 6
Set-TextCell 'D6' '0.9987'
Set-TextCell 'E6' '  -0.10%  '
# Row 7
Set-TextCell 'D7' '0.5014'
Set-TextCell 'E7' '  +1.97%  '
# Row 8
Set-TextCell 'B8' 'Cardano'
Set-TextCell 'C8' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 'D8' '0.2878'
Set-TextCell 'E8' '  +6.85%  '
# Row 9
Set-TextCell 'B9' 'Dogecoin'
Set-TextCell 'C9' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 'D9' '0.06580'
Set-TextCell 'E9' '  +4.56%  '
# Row 10
Set-TextCell 'B10' 'WrappedEther'
Set-TextCell 'C10' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D10' '1.884.02'
Set-TextCell 'E10' '  +5.46%  '
# Row 11
Set-TextCell 'B11' 'Solana'
Set-TextCell 'C11' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 'D11' '17.32'
Set-TextCell 'E11' '  +4.85%  '
# Row 12
Set-TextCell 'B12' 'TRON'
Set-TextCell 'C12' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 'D12' '0.07251'
Set-TextCell 'E12' '  +2.90%  '
# Row 13
Set-TextCell 'B13' 'Polygon'
Set-TextCell 'C13' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D13' '0.6722'
Set-TextCell 'E13' '  +6.93%  '
# Row 14
Set-TextCell 'B14' 'Litecoin'
Set-TextCell 'C14' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D14' '85.20'
Set-TextCell 'E14' '  +6.45%  '
# Row 15
Set-TextCell 'B15' 'Polkadot'
Set-TextCell 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D15' '4.841'
Set-TextCell 'E15' '  +3.84%  '
# Row 16
Set-TextCell 'B16' 'WrappedBTC'
Set-TextCell 'C16' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 'D16' '30.104.43'
Set-TextCell 'E16' '  +6.85%  '
# Row 17
Set-TextCell 'B17' 'Dai'
Set-TextCell 'C17' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D17' '0.9988'
Set-TextCell 'E17' '  -0.13%  '
# Row 18
Set-TextCell 'B18' 'Avalanche'
Set-TextCell 'C18' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 'D18' '12.94'
Set-TextCell 'E18' '  +7.37%  '
# Row 19
Set-TextCell 'B19' 'ShibaInu'
Set-TextCell 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D19' '0.000007560'
Set-TextCell 'E19' '  +4.31%  '
# Row 20
Set-TextCell 'B20' 'BinanceUSD'
Set-TextCell 'C20' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D20' '0.9985'
Set-TextCell 'E20' '  -0.13%  '
# Row 21
Set-TextCell 'B21' 'WrappedliquidstakedEther2.0'
Set-TextCell 'C21' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 'D21' '2.127.82'
Set-TextCell 'E21' '  +5.69%  '
# Row 22
Set-TextCell 'B22' 'Uniswap'
Set-TextCell 'C22' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D22' '4.785'
Set-TextCell 'E22' '  +5.10%  '
# Row 23
Set-TextCell 'B23' 'Chainlink'
Set-TextCell 'C23' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 'D23' '5.560'
Set-TextCell 'E23' '  +5.71%  '
# Row 24
Set-TextCell 'B24' 'Cosmos'
Set-TextCell 'C24' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D24' '9.071'
Set-TextCell 'E24' '  +3.49%  '
# Row 25
Set-TextCell 'B25' 'Monero'
Set-TextCell 'C25' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D25' '145.62'
Set-TextCell 'E25' '  +2.72%  '
# Row 26
Set-TextCell 'B26' 'BitcoinCash'
Set-TextCell 'C26' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 'D26' '135.88'
Set-TextCell 'E26' '  +23.78%  '
# Row 27
Set-TextCell 'B27' 'EthereumClassic'
Set-TextCell 'C27' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D27' '16.84'
Set-TextCell 'E27' '  +6.68%  '
# Row 28
Set-TextCell 'B28' 'LidoDAOToken'
Set-TextCell 'C28' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 'D28' '1.957'
Set-TextCell 'E28' '  +5.44%  '
# Row 29
Set-TextCell 'B29' 'Toncoin'
Set-TextCell 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D29' '1.372'
Set-TextCell 'E29' '  -1.00%  '
# Row 30
Set-TextCell 'B30' 'InternetComputer(DFINITY)'
Set-TextCell 'C30' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D30' '4.209'
Set-TextCell 'E30' '  +0.75%  '
# Row 31
Set-TextCell 'B31' 'Stellar'
Set-TextCell 'C31' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D31' '0.08686'
Set-TextCell 'E31' '  +4.80%  '
# Row 32
Set-TextCell 'B32' 'Filecoin'
Set-TextCell 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D32' '3.958'
Set-TextCell 'E32' '  +4.98%  '
# Row 33
Set-TextCell 'B33' 'Hedera'
Set-TextCell 'C33' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D33' '0.05039'
Set-TextCell 'E33' '  +3.14%  '
# Row 34
Set-TextCell 'B34' 'ARBITRUM'
Set-TextCell 'C34' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D34' '1.145'
Set-TextCell 'E34' '  +5.73%  '
# Row 35
Set-TextCell 'B35' 'ImmutableX'
Set-TextCell 'C35' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D35' '0.6976'
Set-TextCell 'E35' '  +6.42%  '
# Row 36
Set-TextCell 'B36' 'HuobiToken'
Set-TextCell 'C36' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D36' '2.683'
Set-TextCell 'E36' '  +2.52%  '
# Row 37
Set-TextCell 'B37' 'RenderToken'
Set-TextCell 'C37' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D37' '2.295'
Set-TextCell 'E37' '  +10.94%  '
# Row 38
Set-TextCell 'B38' 'MXToken'
Set-TextCell 'C38' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D38' '2.769'
Set-TextCell 'E38' '  +6.05%  '
# Row 39
Set-TextCell 'B39' 'TrustWalletToken'
Set-TextCell 'C39' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D39' '0.9639'
Set-TextCell 'E39' '  +1.77%  '
# Row 40
Set-TextCell 'B40' 'VeChain'
Set-TextCell 'C40' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D40' '0.01637'
Set-TextCell 'E40' '  +5.74%  '
# Row 41
Set-TextCell 'B41' 'FraxShare'
Set-TextCell 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D41' '6.071'
Set-TextCell 'E41' '  +2.47%  '
# Row 42
Set-TextCell 'B42' 'Quant'
Set-TextCell 'C42' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D42' '105.11'
Set-TextCell 'E42' '  +5.30%  '
# Row 43
Set-TextCell 'B43' 'PaxDollar'
Set-TextCell 'C43' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D43' '0.9988'
Set-TextCell 'E43' '  -0.06%  '
# Row 44
Set-TextCell 'D44' '0.4234'
Set-TextCell 'E44' '  +5.96%  '
# Row 45
Set-TextCell 'B45' 'Aptos'
Set-TextCell 'C45' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D45' '7.486'
Set-TextCell 'E45' '  +3.99%  '
# Row 46
Set-TextCell 'B46' 'Algorand'
Set-TextCell 'C46' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D46' '0.1261'
Set-TextCell 'E46' '  +3.72%  '
# Row 47
Set-TextCell 'B47' 'Cronos'
Set-TextCell 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D47' '0.05663'
Set-TextCell 'E47' '  +3.96%  '
# Row 48
Set-TextCell 'B48' 'Elrond'
Set-TextCell 'C48' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 'D48' '32.60'
Set-TextCell 'E48' '  +6.00%  '
# Row 49
Set-TextCell 'D49' '8.313'
Set-TextCell 'E49' '  +3.84%  '
# Row 50
Set-TextCell 'B50' 'Decentraland'
Set-TextCell 'C50' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 'D50' '0.3736'
Set-TextCell 'E50' '  +6.87%  '
# Row 51
Set-TextCell 'B51' 'NEARProtocol'
Set-TextCell 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D51' '1.339'
Set-TextCell 'E51' '  +3.32%  '
